# Add a new "日期：2017.9.27" plan/report block (rows 86-93) to sheet1,
# mirroring the structure of the existing "日期：2017.9.20 四 周三 下午"
# block (rows 49-56), and point the sheet view at the new bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Clone the look & feel (fonts/alignment/number formats) of rows 49:56
#    into the new rows 86:93 so the new block visually matches the other
#    "日期" sections already present in the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A49:D56").Copy()
$ws.Range("A86").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Row heights: header/footer bands match the template block (49:56 ->
# 86:93), while the five data rows need the taller 45pt height because the
# new task description text wraps onto more lines than the template's did.
$ws.Rows.Item(86).RowHeight = $ws.Rows.Item(49).RowHeight
$ws.Rows.Item(87).RowHeight = $ws.Rows.Item(50).RowHeight
$ws.Rows.Item(88).RowHeight = 45
$ws.Rows.Item(89).RowHeight = 45
$ws.Rows.Item(90).RowHeight = 45
$ws.Rows.Item(91).RowHeight = 45
$ws.Rows.Item(92).RowHeight = 45
$ws.Rows.Item(93).RowHeight = $ws.Rows.Item(56).RowHeight

# Re-create the merged header/footer bands for the new block.
$ws.Range("A86:D86").Merge()
$ws.Range("A93:D93").Merge()

# ---------------------------------------------------------------------------
# 2. Fill in the values for the new block.
# ---------------------------------------------------------------------------

# Section date header.
$ws.Range("A86").Value = "日期：2017.9.27 五 周三 上午"

# Column headers.
$ws.Range("A87").Value = "人员"
$ws.Range("B87").Value = "计划任务"
$ws.Range("C87").Value = "完成情况"
$ws.Range("D87").Value = "备注"

# Per-person rows.
$ws.Range("A88").Value = "钟崇尧"
$ws.Range("B88").Value = "学习hbilder开发软件，复习HTML5，CSS3，JavaScript。"

$ws.Range("A89").Value = "黄宝怡"
$ws.Range("B89").Value = "学习hbilder开发软件，复习HTML5，CSS3，JavaScript。"

$ws.Range("A90").Value = "魏仲凯"
$ws.Range("B90").Value = "学习hbilder开发软件，复习HTML5，CSS3，JavaScript。"

$ws.Range("A91").Value = "薛洁鹏"
$ws.Range("B91").Value = "学习hbilder开发软件，复习HTML5，CSS3，JavaScript。"

$ws.Range("A92").Value = "王晓宇"
$ws.Range("B92").Value = "学习hbilder开发软件，复习HTML5，CSS3，JavaScript。"

# Closing "总结：" (summary) row.
$ws.Range("A93").Value = "总结："

# C54 was retyped by the author (same text, "已完成") which is why it now
# points at a distinct shared-string entry upstream; re-apply the value here.
$ws.Range("C54").Value = "已完成"

# ---------------------------------------------------------------------------
# 3. Point the sheet view at the newly added rows, like the author did when
#    they finished the edit.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 73
$ws.Range("C91").Select()
